$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells receiving a numeric-looking replacement value must be pre-formatted as
# Text so Excel stores the new value as a string (matching the source inline
# string cells) instead of silently converting it to a number.
$textCells = @("D5", "D6", "D7", "D9", "D10", "D11", "D12", "D13", "D16", "D17", "D21", "D22", "D23", "D26", "D29", "D31", "D33", "D34", "D35", "D36", "D38", "D40", "D42", "D44", "D45", "D47", "D49", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "42.773.86"
$ws.Range("E2").Value = "  +4.79%  "
$ws.Range("D3").Value = "2.253.59"
$ws.Range("E3").Value = "  +4.45%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "248.92"
$ws.Range("E5").Value = "  +1.67%  "
$ws.Range("D6").Value = "0.634"
$ws.Range("E6").Value = "  +3.41%  "
$ws.Range("D7").Value = "70.72"
$ws.Range("E7").Value = "  +7.95%  "
$ws.Range("E8").Value = "  -0.18%  "
$ws.Range("D9").Value = "0.655"
$ws.Range("E9").Value = "  +17.00%  "
$ws.Range("D10").Value = "39.31"
$ws.Range("E10").Value = "  +11.72%  "
$ws.Range("D11").Value = "59.54"
$ws.Range("E11").Value = "  +3.25%  "
$ws.Range("D12").Value = "0.0966"
$ws.Range("E12").Value = "  +5.75%  "
$ws.Range("D13").Value = "7.44"
$ws.Range("E13").Value = "  +9.37%  "
$ws.Range("E14").Value = "  +0.63%  "
$ws.Range("D15").Value = "2.584.04"
$ws.Range("E15").Value = "  +3.92%  "
$ws.Range("D16").Value = "14.88"
$ws.Range("E16").Value = "  +5.43%  "
$ws.Range("D17").Value = "0.882"
$ws.Range("E17").Value = "  +4.02%  "
$ws.Range("D18").Value = "2.256.26"
$ws.Range("E18").Value = "  +3.76%  "
$ws.Range("D19").Value = "42.705.52"
$ws.Range("E19").Value = "  +4.84%  "
$ws.Range("D20").Value = "0.0₃0990"
$ws.Range("E20").Value = "  +6.66%  "
$ws.Range("D21").Value = "6.30"
$ws.Range("E21").Value = "  +4.38%  "
$ws.Range("D22").Value = "72.97"
$ws.Range("E22").Value = "  +3.09%  "
$ws.Range("D23").Value = "235.46"
$ws.Range("E23").Value = "  +3.52%  "
$ws.Range("E24").Value = "  -1.93%  "
$ws.Range("E25").Value = "  +6.99%  "
$ws.Range("D26").Value = "11.55"
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("E28").Value = "  +1.75%  "
$ws.Range("D29").Value = "3.66"
$ws.Range("E29").Value = "  -1.40%  "
$ws.Range("E30").Value = "  -0.69%  "
$ws.Range("D31").Value = "167.91"
$ws.Range("E31").Value = "  +0.17%  "
$ws.Range("E32").Value = "  +4.45%  "
$ws.Range("D33").Value = "6.48"
$ws.Range("E33").Value = "  +17.77%  "
$ws.Range("D34").Value = "0.125"
$ws.Range("E34").Value = "  +5.76%  "
$ws.Range("D35").Value = "0.0797"
$ws.Range("E35").Value = "  +9.60%  "
$ws.Range("D36").Value = "31.32"
$ws.Range("E36").Value = "  +25.57%  "
$ws.Range("E37").Value = "  +4.46%  "
$ws.Range("D38").Value = "4.42"
$ws.Range("E38").Value = "  +12.19%  "
$ws.Range("E39").Value = "  +4.25%  "
$ws.Range("D40").Value = "0.0323"
$ws.Range("E40").Value = "  +9.56%  "
$ws.Range("E41").Value = "  +7.20%  "
$ws.Range("D42").Value = "12.41"
$ws.Range("E42").Value = "  +9.80%  "
$ws.Range("E43").Value = "  +7.24%  "
$ws.Range("D44").Value = "62.15"
$ws.Range("E44").Value = "  +3.86%  "
$ws.Range("D45").Value = "0.202"
$ws.Range("E45").Value = "  +7.23%  "
$ws.Range("E46").Value = "  +7.37%  "
$ws.Range("D47").Value = "4.87"
$ws.Range("E47").Value = "  +1.18%  "
$ws.Range("E48").Value = "  +4.20%  "
$ws.Range("D49").Value = "1.00"
$ws.Range("E49").Value = "  -0.22%  "
$ws.Range("D50").Value = "1.17"
$ws.Range("E50").Value = "  +3.04%  "
$ws.Range("E51").Value = "  +5.19%  "
